$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the numeric-looking columns (D = Price, E = Volume(1h), G = Hora)
# so Excel stores the new values as literal text, matching the original inline-string cells,
# instead of auto-converting them to numbers/percentages.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "306.84"
$ws.Range("E2").Value = "0.21%"
$ws.Range("F2").Value = "24-1-2023"
$ws.Range("G2").Value = "1"

# Row 3
$ws.Range("D3").Value = "36.32"
$ws.Range("E3").Value = "-0.09%"
$ws.Range("F3").Value = "24-1-2023"
$ws.Range("G3").Value = "1"

# Row 4
$ws.Range("D4").Value = "5.084"
$ws.Range("E4").Value = "0.95%"
$ws.Range("F4").Value = "24-1-2023"
$ws.Range("G4").Value = "1"

# Row 5
$ws.Range("D5").Value = "0.08042"
$ws.Range("E5").Value = "3.08%"
$ws.Range("F5").Value = "24-1-2023"
$ws.Range("G5").Value = "1"

# Row 6
$ws.Range("D6").Value = "2.202"
$ws.Range("E6").Value = "3.70%"
$ws.Range("F6").Value = "24-1-2023"
$ws.Range("G6").Value = "1"

# Row 7
$ws.Range("D7").Value = "8.007"
$ws.Range("E7").Value = "0.79%"
$ws.Range("F7").Value = "24-1-2023"
$ws.Range("G7").Value = "1"

# Row 8
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.156"
$ws.Range("E8").Value = "2.74%"
$ws.Range("F8").Value = "24-1-2023"
$ws.Range("G8").Value = "1"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9280"
$ws.Range("E9").Value = "0.32%"
$ws.Range("F9").Value = "24-1-2023"
$ws.Range("G9").Value = "1"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.09896"
$ws.Range("E10").Value = "1.29%"
$ws.Range("F10").Value = "24-1-2023"
$ws.Range("G10").Value = "1"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1875"
$ws.Range("E11").Value = "-0.01%"
$ws.Range("F11").Value = "24-1-2023"
$ws.Range("G11").Value = "1"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09131"
$ws.Range("E12").Value = "4.82%"
$ws.Range("F12").Value = "24-1-2023"
$ws.Range("G12").Value = "1"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03596"
$ws.Range("E13").Value = "3.07%"
$ws.Range("F13").Value = "24-1-2023"
$ws.Range("G13").Value = "1"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09910"
$ws.Range("E14").Value = "-0.19%"
$ws.Range("F14").Value = "24-1-2023"
$ws.Range("G14").Value = "1"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001435"
$ws.Range("E15").Value = "-0.77%"
$ws.Range("F15").Value = "24-1-2023"
$ws.Range("G15").Value = "1"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005659"
$ws.Range("E16").Value = "0.53%"
$ws.Range("F16").Value = "24-1-2023"
$ws.Range("G16").Value = "1"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.456"
$ws.Range("E17").Value = "-0.61%"
$ws.Range("F17").Value = "24-1-2023"
$ws.Range("G17").Value = "1"

# Row 18
$ws.Range("D18").Value = "2.797"
$ws.Range("E18").Value = "14.68%"
$ws.Range("F18").Value = "24-1-2023"
$ws.Range("G18").Value = "1"

# Row 19
$ws.Range("E19").Value = "-1.29%"
$ws.Range("F19").Value = "24-1-2023"
$ws.Range("G19").Value = "1"

# Row 20
$ws.Range("D20").Value = "0.1319"
$ws.Range("E20").Value = "0.73%"
$ws.Range("F20").Value = "24-1-2023"
$ws.Range("G20").Value = "1"

# Row 21
$ws.Range("D21").Value = "5.064"
$ws.Range("E21").Value = "6.27%"
$ws.Range("F21").Value = "24-1-2023"
$ws.Range("G21").Value = "1"

# Row 22
$ws.Range("D22").Value = "0.2336"
$ws.Range("E22").Value = "2.93%"
$ws.Range("F22").Value = "24-1-2023"
$ws.Range("G22").Value = "1"

# Row 23
$ws.Range("D23").Value = "0.04586"
$ws.Range("E23").Value = "-0.17%"
$ws.Range("F23").Value = "24-1-2023"
$ws.Range("G23").Value = "1"

# Row 24
$ws.Range("E24").Value = "0.78%"
$ws.Range("F24").Value = "24-1-2023"
$ws.Range("G24").Value = "1"

# Row 25
$ws.Range("D25").Value = "0.004753"
$ws.Range("E25").Value = "-6.81%"
$ws.Range("F25").Value = "24-1-2023"
$ws.Range("G25").Value = "1"

# Row 26
$ws.Range("D26").Value = "0.0001299"
$ws.Range("E26").Value = "-7.10%"
$ws.Range("F26").Value = "24-1-2023"
$ws.Range("G26").Value = "1"

# Row 27
$ws.Range("D27").Value = "0.0004498"
$ws.Range("E27").Value = "65.02%"
$ws.Range("F27").Value = "24-1-2023"
$ws.Range("G27").Value = "1"

# Row 28
$ws.Range("F28").Value = "24-1-2023"
$ws.Range("G28").Value = "1"

# Row 29
$ws.Range("F29").Value = "24-1-2023"
$ws.Range("G29").Value = "1"

# Row 30
$ws.Range("F30").Value = "24-1-2023"
$ws.Range("G30").Value = "1"

# Row 31
$ws.Range("F31").Value = "24-1-2023"
$ws.Range("G31").Value = "1"

# Row 32
$ws.Range("F32").Value = "24-1-2023"
$ws.Range("G32").Value = "1"

# Row 33
$ws.Range("F33").Value = "24-1-2023"
$ws.Range("G33").Value = "1"

# Row 34
$ws.Range("F34").Value = "24-1-2023"
$ws.Range("G34").Value = "1"

# Row 35
$ws.Range("F35").Value = "24-1-2023"
$ws.Range("G35").Value = "1"

# Row 36
$ws.Range("F36").Value = "24-1-2023"
$ws.Range("G36").Value = "1"

# Row 37
$ws.Range("F37").Value = "24-1-2023"
$ws.Range("G37").Value = "1"

# Row 38
$ws.Range("F38").Value = "24-1-2023"
$ws.Range("G38").Value = "1"

# Row 39
$ws.Range("D39").Value = "0.01944"
$ws.Range("E39").Value = "6.94%"
$ws.Range("F39").Value = "24-1-2023"
$ws.Range("G39").Value = "1"

# Row 40
$ws.Range("D40").Value = "0.04935"
$ws.Range("E40").Value = "4.08%"
$ws.Range("F40").Value = "24-1-2023"
$ws.Range("G40").Value = "1"

# Row 41
$ws.Range("D41").Value = "0.007763"
$ws.Range("E41").Value = "3.85%"
$ws.Range("F41").Value = "24-1-2023"
$ws.Range("G41").Value = "1"

# Row 42
$ws.Range("E42").Value = "-0.30%"
$ws.Range("F42").Value = "24-1-2023"
$ws.Range("G42").Value = "1"

# Row 43
$ws.Range("D43").Value = "0.007810"
$ws.Range("E43").Value = "1.17%"
$ws.Range("F43").Value = "24-1-2023"
$ws.Range("G43").Value = "1"

# Row 44
$ws.Range("D44").Value = "0.002103"
$ws.Range("E44").Value = "-6.02%"
$ws.Range("F44").Value = "24-1-2023"
$ws.Range("G44").Value = "1"

# Row 45
$ws.Range("D45").Value = "0.01131"
$ws.Range("E45").Value = "7.99%"
$ws.Range("F45").Value = "24-1-2023"
$ws.Range("G45").Value = "1"

# Row 46
$ws.Range("D46").Value = "0.00006240"
$ws.Range("E46").Value = "0.71%"
$ws.Range("F46").Value = "24-1-2023"
$ws.Range("G46").Value = "1"

# Row 47
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "0.04%"
$ws.Range("F47").Value = "24-1-2023"
$ws.Range("G47").Value = "1"

# Row 48
$ws.Range("D48").Value = "29.03"
$ws.Range("E48").Value = "-25.84%"
$ws.Range("F48").Value = "24-1-2023"
$ws.Range("G48").Value = "1"

# Row 49
$ws.Range("D49").Value = "0.001799"
$ws.Range("E49").Value = "-9.94%"
$ws.Range("F49").Value = "24-1-2023"
$ws.Range("G49").Value = "1"

# Row 50
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "0.04%"
$ws.Range("F50").Value = "24-1-2023"
$ws.Range("G50").Value = "1"

# Row 51
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "0.04%"
$ws.Range("F51").Value = "24-1-2023"
$ws.Range("G51").Value = "1"
